$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting (bold, centered, bordered) from H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-16
$values = @(
    @(9, 10),
    @(6, 7),
    @(6, 7),
    @(13, 15),
    @(9, 9),
    @(4, 7),
    @(8, 9),
    @(7, 9),
    @(4, 7),
    @(4, 7),
    @(1, 5),
    @(1, 4),
    @(6, 8),
    @(7, 8),
    @(1, 3)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
